$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress values
$ws.Range("E2").Value = 0.6
$ws.Range("E8").Value = 0.6
$ws.Range("E9").Value = 0.6

# Fill D9 with same value as C9 (the session date), preserving its
# text type (avoid auto date-conversion) by pasting values only
$ws.Range("C9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Update selection to match the target state
$ws.Range("B9:C9").Select()
